# Populate the "Personal Productivity System" workbook:
#   MainMenu  - title, instructions, and hyperlinks to the other sheets
#   Tasks     - header row + one sample task row (with a formatted due date)
#   Routines  - header row
#   Logs      - header row

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# MainMenu sheet
# ---------------------------------------------------------------------------
$wsMenu = $wb.Worksheets.Item("MainMenu")

$wsMenu.Range("A1").Value = "Personal Productivity System"
$wsMenu.Range("A2").Value = "Use the links below to navigate:"
$wsMenu.Range("A4").Value = "Tasks"
$wsMenu.Range("A5").Value = "Routines"
$wsMenu.Range("A6").Value = "Logs"

$wsMenu.Hyperlinks.Add($wsMenu.Range("A4"), "#Tasks!A1", "", "", "Tasks")
$wsMenu.Hyperlinks.Add($wsMenu.Range("A5"), "#Routines!A1", "", "", "Routines")
$wsMenu.Hyperlinks.Add($wsMenu.Range("A6"), "#Logs!A1", "", "", "Logs")

# ---------------------------------------------------------------------------
# Tasks sheet
# ---------------------------------------------------------------------------
$wsTasks = $wb.Worksheets.Item("Tasks")

$wsTasks.Range("A1").Value = "Task ID"
$wsTasks.Range("B1").Value = "Task Name"
$wsTasks.Range("C1").Value = "Due Date"
$wsTasks.Range("D1").Value = "Status"
$wsTasks.Range("E1").Value = "Recurrence"
$wsTasks.Range("F1").Value = "Notes"

$wsTasks.Range("A2").Value = 1
$wsTasks.Range("B2").Value = "Back up budget files"
$wsTasks.Range("C2").Value = 45900
$wsTasks.Range("C2").NumberFormat = "yyyy-mm-dd h:mm:ss"
$wsTasks.Range("D2").Value = "Pending"
$wsTasks.Range("E2").Value = "Monthly"
$wsTasks.Range("F2").Value = ""

# ---------------------------------------------------------------------------
# Routines sheet
# ---------------------------------------------------------------------------
$wsRoutines = $wb.Worksheets.Item("Routines")

$wsRoutines.Range("A1").Value = "Routine Name"
$wsRoutines.Range("B1").Value = "Frequency"
$wsRoutines.Range("C1").Value = "Description"

# ---------------------------------------------------------------------------
# Logs sheet
# ---------------------------------------------------------------------------
$wsLogs = $wb.Worksheets.Item("Logs")

$wsLogs.Range("A1").Value = "Task ID"
$wsLogs.Range("B1").Value = "Task Name"
$wsLogs.Range("C1").Value = "Completed Date"
$wsLogs.Range("D1").Value = "Original Due"
$wsLogs.Range("E1").Value = "Notes"
